# Update recalculated market-price / profit figures across the leve
# profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), per the
# scheduled runner's refreshed pricing snapshot.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 46.625
$ws.Range("I6").Value = 46.625
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 139.875
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = ""
$ws.Range("N6").Value = -27.875

$ws.Range("H19").Value = 203
$ws.Range("I19").Value = 265.8
$ws.Range("J19").Value = 98.333336
$ws.Range("K19").Value = 265.8
$ws.Range("L19").Value = 98.333336
$ws.Range("M19").Value = -90.80000000000001
$ws.Range("N19").Value = -448.333336

$ws.Range("H33").Value = 624.5
$ws.Range("I33").Value = 399.8
$ws.Range("K33").Value = 399.8
$ws.Range("M33").Value = -170.8

$ws.Range("H51").Value = 6000
$ws.Range("J51").Value = 6000
$ws.Range("L51").Value = 6000
$ws.Range("N51").Value = -6968

$ws.Range("H112").Value = 2659.9

$ws.Range("H135").Value = 1291.9231
$ws.Range("I135").Value = 1145
$ws.Range("J135").Value = 2100
$ws.Range("K135").Value = 10305
$ws.Range("L135").Value = 18900
$ws.Range("M135").Value = -7770
$ws.Range("N135").Value = -23970

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 7271.75
$ws.Range("I3").Value = 7271.75
$ws.Range("K3").Value = 7271.75
$ws.Range("M3").Value = -7156.75

$ws.Range("H16").Value = 14625.444
$ws.Range("I16").Value = 16938.166
$ws.Range("J16").Value = 10000
$ws.Range("K16").Value = 16938.166
$ws.Range("L16").Value = 10000
$ws.Range("M16").Value = -16651.166
$ws.Range("N16").Value = -10574

$ws.Range("H32").Value = 5889.8887
$ws.Range("I32").Value = 4629.6
$ws.Range("K32").Value = 4629.6
$ws.Range("M32").Value = -4342.6

$ws.Range("H132").Value = 1411.7858
$ws.Range("I132").Value = 1356.1666
$ws.Range("J132").Value = 1745.5
$ws.Range("K132").Value = 4068.4998
$ws.Range("L132").Value = 5236.5
$ws.Range("M132").Value = -1538.4998
$ws.Range("N132").Value = -10296.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 765
$ws.Range("I10").Value = 647.5
$ws.Range("J10").Value = 1000
$ws.Range("K10").Value = 647.5
$ws.Range("L10").Value = 1000
$ws.Range("M10").Value = -507.5
$ws.Range("N10").Value = -1280

$ws.Range("H20").Value = 1360.5385
$ws.Range("I20").Value = 680.8889
$ws.Range("J20").Value = 2889.75
$ws.Range("K20").Value = 680.8889
$ws.Range("L20").Value = 2889.75
$ws.Range("M20").Value = -433.8889
$ws.Range("N20").Value = -3383.75

$ws.Range("H135").Value = 48499.5
$ws.Range("J135").Value = 48499.5
$ws.Range("L135").Value = 48499.5
$ws.Range("N135").Value = -58639.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 181
$ws.Range("I10").Value = 181
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 181
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = ""
$ws.Range("N10").Value = -42

$ws.Range("H15").Value = 40739
$ws.Range("I15").Value = 40739
$ws.Range("K15").Value = 40739
$ws.Range("M15").Value = -40569

$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").Value = ""

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").Value = ""

$ws.Range("H134").Value = 2377
$ws.Range("I134").Value = 1323.2307
$ws.Range("J134").Value = 4660.1665
$ws.Range("K134").Value = 3969.6921
$ws.Range("L134").Value = 13980.4995
$ws.Range("M134").Value = -1434.6921
$ws.Range("N134").Value = -19050.4995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 395
$ws.Range("I12").Value = 226.6
$ws.Range("J12").Value = 563.4
$ws.Range("K12").Value = 679.8
$ws.Range("L12").Value = 1690.2
$ws.Range("M12").Value = -506.8
$ws.Range("N12").Value = -2036.2

$ws.Range("H132").Value = 2833.3333
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 286.66666
$ws.Range("I13").Value = 55
$ws.Range("J13").Value = 750
$ws.Range("K13").Value = 55
$ws.Range("L13").Value = 750
$ws.Range("M13").Value = 84
$ws.Range("N13").Value = -1028

$ws.Range("H26").Value = 30042
$ws.Range("J26").Value = 30042
$ws.Range("L26").Value = 30042
$ws.Range("N26").Value = -30602

$ws.Range("H50").Value = 30042
$ws.Range("J50").Value = 30042
$ws.Range("L50").Value = 30042
$ws.Range("N50").Value = -31038

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 16107
$ws.Range("I43").Value = 9000
$ws.Range("J43").Value = 17528.4
$ws.Range("K43").Value = 9000
$ws.Range("L43").Value = 17528.4
$ws.Range("M43").Value = -8807
$ws.Range("N43").Value = -17914.4

$ws.Range("H46").Value = 3583.1667
$ws.Range("I46").Value = 3000
$ws.Range("J46").Value = 3777.5557
$ws.Range("K46").Value = 3000
$ws.Range("L46").Value = 3777.5557
$ws.Range("M46").Value = -2812
$ws.Range("N46").Value = -4153.5557

$ws.Range("H55").Value = 1625.9
$ws.Range("I55").Value = 1004.1667
$ws.Range("J55").Value = 2558.5
$ws.Range("K55").Value = 1004.1667
$ws.Range("L55").Value = 2558.5
$ws.Range("M55").Value = -831.1667
$ws.Range("N55").Value = -2904.5

$ws.Range("H68").Value = 3950
$ws.Range("I68").Value = 3950
$ws.Range("K68").Value = 3950
$ws.Range("M68").Value = -3201

$ws.Range("H71").Value = 3950
$ws.Range("I71").Value = 3950
$ws.Range("K71").Value = 19750
$ws.Range("M71").Value = -16006

$ws.Range("H122").Value = 3196.4736
$ws.Range("I122").Value = 3114.8462
$ws.Range("K122").Value = 9344.5386
$ws.Range("M122").Value = -6894.5386

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 15958.7
$ws.Range("I41").Value = 16218.6
$ws.Range("J41").Value = 15698.8
$ws.Range("K41").Value = 16218.6
$ws.Range("L41").Value = 15698.8
$ws.Range("M41").Value = -15828.6
$ws.Range("N41").Value = -16478.8

$ws.Range("H126").Value = 3002.476
$ws.Range("I126").Value = 3018.3076
$ws.Range("J126").Value = 2976.75
$ws.Range("K126").Value = 9054.9228
$ws.Range("L126").Value = 8930.25
$ws.Range("M126").Value = -6584.9228
$ws.Range("N126").Value = -13870.25
